$wb = $excel.ActiveWorkbook

# --- BOM_PSC sheet: insert new component rows (4-7) ---
$ws1 = $wb.Worksheets.Item("BOM_PSC")

$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "Phoenix 1840502"
$ws1.Range("C4").Value = "16-pin 3.5mm terminal block plug"

$ws1.Range("A5").Value = 8
$ws1.Range("C5").Value = "26AWG twisted pair wire"

$ws1.Range("A6").Value = 2
$ws1.Range("B6").Value = "Phoenix 1840366"
$ws1.Range("C6").Value = "2-pin 3.5mm terminal block plug"

$ws1.Range("A7").Value = 1
$ws1.Range("C7").Value = "26AWG twisted pair wire"

# --- Update selections on both sheets ---
$ws2 = $wb.Worksheets.Item("BOM_PWMAttenuation")
$ws2.Range("F17").Select() | Out-Null

$ws1.Range("B11").Select() | Out-Null
